$pres = $ppt.Presentations
Write-Output ($pres | Get-Member | Out-String)
